$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Price (column D) updates - set as text to preserve exact formatting
Set-TextValue 'D2' '68.730.10'
Set-TextValue 'D3' '3.746.60'
Set-TextValue 'D5' '601.60'
Set-TextValue 'D6' '169.16'
Set-TextValue 'D7' '3.746.63'
Set-TextValue 'D13' '38.32'
Set-TextValue 'D14' '0.0000248'
Set-TextValue 'D15' '4.375.07'
Set-TextValue 'D16' '3.753.90'
Set-TextValue 'D17' '68.741.08'
Set-TextValue 'D20' '17.09'
Set-TextValue 'D22' '495.72'
Set-TextValue 'D24' '0.0000154'
Set-TextValue 'D25' '85.33'
Set-TextValue 'D28' '10.34'
Set-TextValue 'D31' '2.52'
Set-TextValue 'D33' '31.85'
Set-TextValue 'D34' '3.894.14'
Set-TextValue 'D36' '3.681.72'
Set-TextValue 'D38' '1.01'
Set-TextValue 'D42' '438.40'
Set-TextValue 'D43' '48.85'
Set-TextValue 'D44' '2.93'
Set-TextValue 'D48' '40.60'
Set-TextValue 'D49' '141.53'
Set-TextValue 'D50' '2.790.80'

# Volume(1h) (column E) updates - plain text assignment (non-numeric strings)
$ws.Range('E2').Value = '  +2.54%  '
$ws.Range('E3').Value = '  +2.42%  '
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('E5').Value = '  +1.70%  '
$ws.Range('E6').Value = '  +1.90%  '
$ws.Range('E7').Value = '  +2.39%  '
$ws.Range('E8').Value = '  -0.11%  '
$ws.Range('E9').Value = '  +2.77%  '
$ws.Range('E10').Value = '  +4.31%  '
$ws.Range('E11').Value = '  +3.41%  '
$ws.Range('E12').Value = '  +0.42%  '
$ws.Range('E13').Value = '  +1.80%  '
$ws.Range('E14').Value = '  +3.88%  '
$ws.Range('E15').Value = '  +2.14%  '
$ws.Range('E16').Value = '  +2.22%  '
$ws.Range('E17').Value = '  +2.39%  '
$ws.Range('E19').Value = '  +0.90%  '
$ws.Range('E20').Value = '  +2.00%  '
$ws.Range('E21').Value = '  +19.80%  '
$ws.Range('E22').Value = '  +0.98%  '
$ws.Range('E23').Value = '  +2.33%  '
$ws.Range('E24').Value = '  +11.26%  '
$ws.Range('E25').Value = '  +0.49%  '
$ws.Range('E26').Value = '  +1.77%  '
$ws.Range('E27').Value = '  +1.80%  '
$ws.Range('E28').Value = '  +4.28%  '
$ws.Range('E29').Value = '  +0.49%  '
$ws.Range('E30').Value = '  +3.02%  '
$ws.Range('E31').Value = '  +7.12%  '
$ws.Range('E32').Value = '  +2.89%  '
$ws.Range('E33').Value = '  +0.77%  '
$ws.Range('E34').Value = '  +2.38%  '
$ws.Range('E35').Value = '  +2.20%  '
$ws.Range('E36').Value = '  +2.25%  '
$ws.Range('E37').Value = '  -0.17%  '
$ws.Range('E38').Value = '  +3.58%  '
$ws.Range('E39').Value = '  +1.94%  '
$ws.Range('E40').Value = '  +1.29%  '
$ws.Range('E41').Value = '  +1.10%  '
$ws.Range('E42').Value = '  -0.33%  '
$ws.Range('E43').Value = '  +0.89%  '
$ws.Range('E44').Value = '  +6.20%  '
$ws.Range('E45').Value = '  +3.17%  '
$ws.Range('E46').Value = '  +2.25%  '
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('E48').Value = '  +2.78%  '
$ws.Range('E49').Value = '  +0.02%  '
$ws.Range('E50').Value = '  +1.65%  '
$ws.Range('E51').Value = '  +2.83%  '
